# "User Stories" workbook - add the two new "worded month" user stories
# and mark the existing one (and the two new ones) as Done.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 34 (existing story #20): it is now finished -----------------------
$ws.Range("G34").Value = "Done"          # status: In-progress -> Done
$ws.Range("I34").Value = 43201           # Date Accomplished

# --- Row 35 (new story #20.1): worded months in time-in tables -------------
$ws.Range("A35").Value = 20.1
$ws.Range("B35").Value = "user"
$ws.Range("C35").Value = "to view numerical months as worded month in time-in tables"
$ws.Range("D35").Value = "so that I can easily distinguish easily view months"
$ws.Range("F35").Value = "Low"
$ws.Range("G35").Value = "Done"
$ws.Range("H35").Value = 43182           # Date Requested
$ws.Range("I35").Value = 43201           # Date Accomplished

# --- Row 36 (new story #20.2): worded months in report generations ---------
$ws.Range("A36").Value = 20.2
$ws.Range("B36").Value = "user"
$ws.Range("C36").Value = "to view numerical months as worded month in report generations"
$ws.Range("D36").Value = "so that I can easily distinguish easily view months"
$ws.Range("F36").Value = "Low"
$ws.Range("G36").Value = "Done"
$ws.Range("H36").Value = 43182           # Date Requested
$ws.Range("I36").Value = 43201           # Date Accomplished

# --- Give the two new rows the same look as the other "Done" rows ----------
# (row 30 is a finished story already styled with the green fill + date
# format used for completed entries)
$ws.Range("A30:I30").Copy()
$ws.Range("A34:I36").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# --- Scroll/selection state so the newly added rows are in view ------------
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 15
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("H40").Select()
